# adicionando nova versão de login e página de comunidade
#
# 1) New weekly-status bullet "Adicionando nova versão de login e tela de
#    comunidade;" + dates, right before the two blank list paragraphs
#    that follow "Página de login;".
# 2) "Futuras pautas:" now starts a new page: three blank (b/bCs, sz 28)
#    paragraphs are inserted before it and the heading run is re-created
#    with a <w:lastRenderedPageBreak/> in front of it.
# 3) New bullet "Desenvolvimento final do site;" appended to the
#    "Futuras pautas" list, right after the HLD/LLD bullet.

$d = $word.ActiveDocument

function Wrap-PkgXml([string]$body) {
    $head = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $tail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $head + $body + $tail
}

# Plain paragraph-index lookup (avoids the Find-collapsed-range/.Next()
# quirks of this COM shim) - returns the 1-based index of the first
# paragraph whose text contains $pattern.
function Get-ParagraphIndex([string]$pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Change 1 (applied first - lowest in the doc, so it doesn't shift the
# anchors used by the edits below it in the script):
# "Desenvolvimento final do site;" bullet, numId 8, right after the
# "Desenvolvimento do HLD / LLD do projeto;" bullet.
# ---------------------------------------------------------------------
$idxHld = Get-ParagraphIndex "Desenvolvimento do HLD / LLD do projeto;"
$idxAfterHld = $idxHld + 1
$d.Paragraphs.Item($idxAfterHld).Range.InsertParagraphBefore()

$bodySite = '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Desenvolvimento final do site;</w:t></w:r></w:p>'
$d.Paragraphs.Item($idxAfterHld).Range.InsertXML((Wrap-PkgXml $bodySite))

# ---------------------------------------------------------------------
# Change 2: "Futuras pautas:" gets pushed onto its own page - strip the
# text run from the existing paragraph, then add three blank paragraphs
# followed by one holding the run again (now with lastRenderedPageBreak).
# ---------------------------------------------------------------------
$headingPPr = '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>'

$idxFuturas = Get-ParagraphIndex "Futuras pautas:"
$emptyHeadingBody = '<w:p>' + $headingPPr + '</w:p>'
$d.Paragraphs.Item($idxFuturas).Range.InsertXML((Wrap-PkgXml $emptyHeadingBody))

$idxAfterFuturas = $idxFuturas + 1
$d.Paragraphs.Item($idxAfterFuturas).Range.InsertParagraphBefore()

$blankHeading = '<w:p>' + $headingPPr + '</w:p>'
$filledHeading = '<w:p>' + $headingPPr + '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>Futuras pautas:</w:t></w:r></w:p>'
$bodyHeading = $blankHeading + $blankHeading + $blankHeading + $filledHeading
$d.Paragraphs.Item($idxAfterFuturas).Range.InsertXML((Wrap-PkgXml $bodyHeading))

# ---------------------------------------------------------------------
# Change 3: new bullet "Adicionando nova versão de login e tela de
# comunidade;" (numId 7) + blank PargrafodaLista paragraph + blank
# ind-left-360 paragraph, inserted right before the two blank
# PargrafodaLista paragraphs that follow "Página de login;".
# ---------------------------------------------------------------------
$idxLogin = Get-ParagraphIndex "P.gina de login"
$idxAfterLogin = $idxLogin + 1
$d.Paragraphs.Item($idxAfterLogin).Range.InsertParagraphBefore()

$bodyLogin = '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Adicionando nova versão de login e tela de comunidade;</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>[23/05]</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$d.Paragraphs.Item($idxAfterLogin).Range.InsertXML((Wrap-PkgXml $bodyLogin))
